$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "NA" value that was previously in C114 (it moves to the new row)
$ws.Range("C114").Value = ""

# Add the new row 115 with the same pattern as the other "nothing relevant" rows
# (leading apostrophe forces text so the date-like string isn't auto-converted
# to a date serial number; ClearFormats drops the resulting quote-prefix style
# so the cell ends up with no explicit style, like its neighbours)
$ws.Range("A115").Value = "'2025-05-21"
$ws.Range("A115").ClearFormats()
$ws.Range("B115").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C115").Value = "NA"
$ws.Range("D115").Value = 1
